$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64; this shifts the existing rows 64-95 down to 65-96,
# matching the rest of the diff (which is just that shift) and growing the used
# range from A1:R95 to A1:R96.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record. The constant columns
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Calidad, Unidad de comercializacion, Kg o Unidades, Clasificacion) match the
# rest of the block exactly.
$ws.Range("A64").Value = 10
$ws.Range("B64").Value = "Vega Modelo de Temuco"
$ws.Range("C64").Value = "La Araucanía"
$ws.Range("D64").Value = 44755
$ws.Range("E64").Value = 9
$ws.Range("F64").Value = 100112035
$ws.Range("G64").Value = "Bruselas (repollito)"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 40
$ws.Range("K64").Value = 26000
$ws.Range("L64").Value = 26000
$ws.Range("M64").Value = 26000
$ws.Range("N64").Value = "$/malla 10 kilos"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 2600
$ws.Range("Q64").Value = 10
$ws.Range("R64").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the column (style "s=2").
$ws.Range("D64").NumberFormat = $ws.Range("D65").NumberFormat
